$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new numeric-looking value would
# otherwise be auto-converted to a real number by Excel, so they stay text
# like the rest of the Price column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '43.708.28'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '2.307.01'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').Value = '94.89'
$ws.Range('E5').Value = '  +8.69%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '268.70'
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '0.622'
$ws.Range('E9').Value = '  +1.49%  '
$ws.Range('D10').Value = '44.81'
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('D11').Value = '0.0932'
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').Value = '8.09'
$ws.Range('E12').Value = '  +6.67%  '
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').Value = '2.654.32'
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('D15').Value = '15.32'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').Value = '0.859'
$ws.Range('E16').Value = '  +6.48%  '
$ws.Range('D17').Value = '2.307.83'
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('D18').Value = '43.625.68'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').Value = '0.0000107'
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('D20').Value = '6.31'
$ws.Range('E20').Value = '  +4.07%  '
$ws.Range('D21').Value = '71.32'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').Value = '238.21'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('D23').Value = '2.28'
$ws.Range('E23').Value = '  -5.25%  '
$ws.Range('D24').Value = '9.50'
$ws.Range('E24').Value = '  +6.62%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '11.27'
$ws.Range('E26').Value = '  +3.09%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = '2.50'
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('D28').Value = '2.33'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '3.40'
$ws.Range('E29').Value = '  -4.81%  '
$ws.Range('D30').Value = '38.60'
$ws.Range('E30').Value = '  -2.99%  '
$ws.Range('D31').Value = '22.58'
$ws.Range('E31').Value = '  +7.32%  '
$ws.Range('D32').Value = '171.60'
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('D33').Value = '0.0895'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('E34').Value = '  +2.59%  '
$ws.Range('E35').Value = '  +1.46%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.109'
$ws.Range('E36').Value = '  -3.08%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.0358'
$ws.Range('E37').Value = '  +1.55%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '4.47'
$ws.Range('E38').Value = '  +1.69%  '
$ws.Range('D39').Value = '3.45'
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '2.32'
$ws.Range('E40').Value = '  +3.52%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.234'
$ws.Range('E41').Value = '  +14.34%  '
$ws.Range('E42').Value = '  +17.06%  '
$ws.Range('D43').Value = '12.10'
$ws.Range('E43').Value = '  -4.48%  '
$ws.Range('D44').Value = '5.46'
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('D45').Value = '61.75'
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('D46').Value = '9.07'
$ws.Range('E46').Value = '  +6.26%  '
$ws.Range('E47').Value = '  +2.78%  '
$ws.Range('D48').Value = '100.30'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('D50').Value = '2.532.96'
$ws.Range('E50').Value = '  +1.73%  '
$ws.Range('D51').Value = '0.425'
$ws.Range('E51').Value = '  -1.83%  '
